# Applies swapped values across paired stock rows (same item, different batch/date)
# as described by the commit diff. Each pair/triple of rows sharing the same
# product description has its Item Code (B), Sale Rate (D, when applicable),
# MRP/Rate (E), Quantity (F) and Value (G) columns rotated among the rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

${ws}.Range("B136").Value = 48654
${ws}.Range("E136").Value = 38.26
${ws}.Range("F136").Value = -1
${ws}.Range("G136").Value = -32.02
${ws}.Range("B137").Value = 63902
${ws}.Range("E137").Value = 34.04
${ws}.Range("F137").Value = 2
${ws}.Range("G137").Value = 64.04000000000001
${ws}.Range("B233").Value = 64979
${ws}.Range("E233").Value = 314.41
${ws}.Range("F233").Value = 26
${ws}.Range("G233").Value = 7689.5
${ws}.Range("B234").Value = 48719
${ws}.Range("E234").Value = 353.35
${ws}.Range("F234").Value = -81
${ws}.Range("G234").Value = -23955.75
${ws}.Range("B246").Value = 48706
${ws}.Range("E246").Value = 39.8
${ws}.Range("F246").Value = -144
${ws}.Range("G246").Value = -4795.2
${ws}.Range("B247").Value = 64973
${ws}.Range("E247").Value = 35.4
${ws}.Range("F247").Value = 104
${ws}.Range("G247").Value = 3463.2
${ws}.Range("B277").Value = 63565
${ws}.Range("E277").Value = 109.19
${ws}.Range("F277").Value = 60
${ws}.Range("G277").Value = 6162.6
${ws}.Range("B278").Value = 61610
${ws}.Range("E278").Value = 122.71
${ws}.Range("F278").Value = -58
${ws}.Range("G278").Value = -5957.18
${ws}.Range("B292").Value = 55373
${ws}.Range("E292").Value = 163.62
${ws}.Range("F292").Value = -94
${ws}.Range("G292").Value = -13562.32
${ws}.Range("B293").Value = 63520
${ws}.Range("E293").Value = 153.4
${ws}.Range("F293").Value = 81
${ws}.Range("G293").Value = 11686.68
${ws}.Range("B294").Value = 63531
${ws}.Range("E294").Value = 152.53
${ws}.Range("F294").Value = 80
${ws}.Range("G294").Value = 11478.4
${ws}.Range("B295").Value = 57802
${ws}.Range("E295").Value = 162.71
${ws}.Range("F295").Value = -79
${ws}.Range("G295").Value = -11334.92
${ws}.Range("B296").Value = 63571
${ws}.Range("F296").Value = 9
${ws}.Range("G296").Value = 1291.32
${ws}.Range("B420").Value = 58047
${ws}.Range("D420").Value = 105.54
${ws}.Range("E420").Value = 126.1
${ws}.Range("F420").Value = 43
${ws}.Range("G420").Value = 4538.22
${ws}.Range("B421").Value = 47097
${ws}.Range("D421").Value = 112.28
${ws}.Range("E421").Value = 134.16
${ws}.Range("F421").Value = 15
${ws}.Range("G421").Value = 1684.2
${ws}.Range("B465").Value = 65069
${ws}.Range("E465").Value = 14.3
${ws}.Range("F465").Value = 2
${ws}.Range("G465").Value = 26.9
${ws}.Range("B466").Value = 53757
${ws}.Range("E466").Value = 16.08
${ws}.Range("F466").Value = -159
${ws}.Range("G466").Value = -2138.55
${ws}.Range("B467").Value = 65068
${ws}.Range("E467").Value = 13.97
${ws}.Range("F467").Value = 128
${ws}.Range("G467").Value = 1683.2
${ws}.Range("B468").Value = 53602
${ws}.Range("E468").Value = 15.69
${ws}.Range("F468").Value = -231
${ws}.Range("G468").Value = -3037.65
${ws}.Range("B472").Value = 64915
${ws}.Range("E472").Value = 20.98
${ws}.Range("F472").Value = 0
${ws}.Range("G472").Value = 0
${ws}.Range("B473").Value = 45695
${ws}.Range("E473").Value = 23.58
${ws}.Range("F473").Value = -36
${ws}.Range("G473").Value = -710.28
${ws}.Range("B479").Value = 45718
${ws}.Range("E479").Value = 19.38
${ws}.Range("F479").Value = -294
${ws}.Range("G479").Value = -4768.68
${ws}.Range("B480").Value = 64927
${ws}.Range("E480").Value = 17.26
${ws}.Range("F480").Value = 222
${ws}.Range("G480").Value = 3600.84
${ws}.Range("B591").Value = 60031
${ws}.Range("E591").Value = 111.69
${ws}.Range("F591").Value = -5
${ws}.Range("G591").Value = -492.5
${ws}.Range("B592").Value = 64836
${ws}.Range("E592").Value = 104.71
${ws}.Range("F592").Value = 3
${ws}.Range("G592").Value = 295.5
${ws}.Range("B596").Value = 64830
${ws}.Range("E596").Value = 34.9
${ws}.Range("F596").Value = 113
${ws}.Range("G596").Value = 3709.79
${ws}.Range("B597").Value = 60022
${ws}.Range("E597").Value = 37.22
${ws}.Range("F597").Value = -113
${ws}.Range("G597").Value = -3709.79
${ws}.Range("B705").Value = 63150
${ws}.Range("D705").Value = 75.68000000000001
${ws}.Range("E705").Value = 80.45
${ws}.Range("F705").Value = 91
${ws}.Range("G705").Value = 6886.88
${ws}.Range("B706").Value = 61428
${ws}.Range("D706").Value = 69.16
${ws}.Range("E706").Value = 73.52
${ws}.Range("F706").Value = 1
${ws}.Range("G706").Value = 69.16
${ws}.Range("B732").Value = 65079
${ws}.Range("F732").Value = 21
${ws}.Range("G732").Value = 858.27
${ws}.Range("B733").Value = 65362
${ws}.Range("F733").Value = 69
${ws}.Range("G733").Value = 2820.03
